$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores figures as literal text (e.g. "309.80", "42.972.18"
# using dots as thousands separators) rather than numbers, so trailing zeros and
# the multi-dot notation survive. Pre-format the numeric-looking cells as Text
# (one Range.NumberFormat call per cell - a comma-joined multi-area address only
# applies to its first area) before writing, so Excel does not silently convert
# them to numeric values.
$textPriceCells = @("D5", "D6", "D7", "D10", "D11", "D13", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D35", "D36", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.972.18"
$ws.Range("E2").Value = "  +2.10%  "

$ws.Range("D3").Value = "2.301.82"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "309.80"
$ws.Range("E5").Value = "  +1.37%  "

$ws.Range("D6").Value = "100.01"
$ws.Range("E6").Value = "  +4.45%  "

$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +5.60%  "

$ws.Range("D10").Value = "36.08"
$ws.Range("E10").Value = "  +2.62%  "

$ws.Range("D11").Value = "0.0818"
$ws.Range("E11").Value = "  +3.56%  "

$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "7.15"
$ws.Range("E13").Value = "  +7.83%  "

$ws.Range("D14").Value = "2.660.22"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "14.86"
$ws.Range("E15").Value = "  +3.39%  "

$ws.Range("D16").Value = "2.304.84"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").Value = "0.800"
$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").Value = "42.938.40"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("E20").Value = "  +1.95%  "

$ws.Range("D21").Value = "6.06"
$ws.Range("E21").Value = "  +1.59%  "

$ws.Range("D22").Value = "68.18"
$ws.Range("E22").Value = "  +0.79%  "

$ws.Range("D23").Value = "239.53"
$ws.Range("E23").Value = "  +0.64%  "

$ws.Range("D24").Value = "2.01"
$ws.Range("E24").Value = "  +3.97%  "

$ws.Range("E25").Value = "  +1.52%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "24.12"
$ws.Range("E27").Value = "  +1.75%  "

$ws.Range("D28").Value = "38.54"
$ws.Range("E28").Value = "  +5.22%  "

$ws.Range("D29").Value = "9.64"
$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("E30").Value = "  +0.28%  "

$ws.Range("D31").Value = "168.22"
$ws.Range("E31").Value = "  +5.16%  "

$ws.Range("D32").Value = "5.33"
$ws.Range("E32").Value = "  +1.96%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("D35").Value = "17.67"
$ws.Range("E35").Value = "  +3.36%  "

$ws.Range("D36").Value = "0.0738"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value = "0.116"
$ws.Range("E39").Value = "  +1.50%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.83"
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("E41").Value = "  +5.16%  "

$ws.Range("D42").Value = "2.30"
$ws.Range("E42").Value = "  -4.65%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.967.79"
$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0288"
$ws.Range("E44").Value = "  +1.75%  "

$ws.Range("D45").Value = "19.14"
$ws.Range("E45").Value = "  +0.45%  "

$ws.Range("D46").Value = "3.01"
$ws.Range("E46").Value = "  +2.64%  "

$ws.Range("D47").Value = "9.83"
$ws.Range("E47").Value = "  -1.02%  "

$ws.Range("E48").Value = "  +18.32%  "

$ws.Range("D49").Value = "55.07"
$ws.Range("E49").Value = "  +3.44%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.530.35"
$ws.Range("E50").Value = "  +1.71%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "1.55"
$ws.Range("E51").Value = "  +2.49%  "

